# This script appends a new "Solutions" section to the end of the document body,
# mirroring the existing "Use Cases for CQL 2.0" section (page break + Heading1,
# followed by two Heading2 sub-sections with Thoughts/Proposal/To-do write-ups).

$d = $word.ActiveDocument

# Create a fresh, empty trailing paragraph (just before the final section break)
# to use as the insertion point for the new content.
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()
$target = $d.Paragraphs.Last

# Build the new body content as a raw WordprocessingML fragment (paragraph by
# paragraph) so headings, tabs and the page break come through as real OOXML
# elements rather than plain-text approximations.
$bodyXml = ""
$bodyXml += '<w:p><w:r><w:br w:type="page"/></w:r></w:p>'
$bodyXml += '<w:p><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>Solutions</w:t></w:r></w:p>'
$bodyXml += '<w:p/>'
$bodyXml += '<w:p><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r><w:t>Associated Object Retrieval</w:t></w:r></w:p>'
$bodyXml += '<w:p/>'
$bodyXml += '<w:p><w:r><w:t>Thoughts:  This is very common request, so it deserves high priority.  I propose adding an optional element to the root of a CQL query (like Query Modifiers are now) to specify which associations are to be returned.  The association population should be configurable to either populate ALL associations up to a certain number of levels, or a named association (or multiple associations).  In the case of named associations, the query developer should be able to specify the names of sub-associations to populate as well, in a recursive fashion.</w:t></w:r></w:p>'
$bodyXml += '<w:p><w:r><w:t>Proposal:  The schema AssociationPopulationSpec.xsd fulfills the requirements.  It allows for a choice of depth-based population or named association population.  The naming is specified recursively, and the depth based population is a simple integer value.  The schema presently allows for a flag to be set indicating “infinite” depth population, the usefulness of which needs to be evaluated.</w:t></w:r></w:p>'
$bodyXml += '<w:p><w:r><w:t>To-do:</w:t></w:r></w:p>'
$bodyXml += '<w:p><w:r><w:tab/><w:t>I. Evaluate the usefulness and practicality of implementation of the ‘infinite depth’ flag</w:t></w:r></w:p>'
$bodyXml += '<w:p><w:r><w:tab/><w:t>II. Naming of elements and types in the schema needs some work</w:t></w:r></w:p>'
$bodyXml += '<w:p/>'
$bodyXml += '<w:p/>'
$bodyXml += '<w:p><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r><w:t>Temporal Queries</w:t></w:r></w:p>'
$bodyXml += '<w:p><w:pPr><w:pStyle w:val="Heading2"/></w:pPr></w:p>'
$bodyXml += '<w:p><w:r><w:t>Thoughts: In the TBPT case, this might be a modeling issue.  Things like “age” should never be stored in a database, since they change as time progresses.  From a query perspective, we would need a way to make a query “relative to” some other value.  This gets into the area of joins, which CQL doesn’t really do.  In this case, the value is relative to today’s date, so it’s a known value and not really a join but a value replacement on the server side.</w:t></w:r></w:p>'
$bodyXml += '<w:p><w:r><w:t>To-do:</w:t></w:r></w:p>'
$bodyXml += '<w:p><w:r><w:tab/><w:t>I. Evaluate some TBPT models to see if anybody actually stores “age” values.</w:t></w:r></w:p>'
$bodyXml += '<w:p><w:r><w:tab/><w:t>II. Develop a specialized query type for temporal queries</w:t></w:r></w:p>'
$bodyXml += '<w:p><w:r><w:tab/></w:r><w:r><w:tab/><w:t>A.  Might be dependent on the strongly typed queries values request</w:t></w:r></w:p>'
$bodyXml += '<w:p/>'

# Wrap the fragment in the Flat-OPC envelope Range.InsertXML requires, then
# insert it in place of the empty target paragraph.
$flatOpc = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' +
           $bodyXml +
           '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$target.Range.InsertXML($flatOpc)
